$d = $word.ActiveDocument

# Hunk 1: sobriety clause - add drug/alcohol testing requirement
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$old1 = "Defendant shall maintain sobriety while on bond, and shall not possess, consume, or purchase alcohol or drugs of abuse."
$new1 = "Defendant shall maintain sobriety while on bond, shall not possess, consume, or purchase alcohol or drugs of abuse, and shall submit to drug and alcohol testing as directed by the Office of Community Control."
$find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
